$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.933.12"
$ws.Range("E2").Value = "  -3.67%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.37"
$ws.Range("E3").Value = "  -2.72%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.50"
$ws.Range("E5").Value = "  -2.30%  "

# Row 6
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("E7").Value = "  -5.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3722"
$ws.Range("E8").Value = "  -2.71%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07474"
$ws.Range("E9").Value = "  -3.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9347"
$ws.Range("E10").Value = "  -4.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.35"
$ws.Range("E11").Value = "  -3.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.83"
$ws.Range("E12").Value = "  -0.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.735"
$ws.Range("E13").Value = "  -3.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.442"
$ws.Range("E14").Value = "  -4.40%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06855"
$ws.Range("E15").Value = "  -1.83%  "

# Row 16
$ws.Range("E16").Value = "  -0.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.61"
$ws.Range("E17").Value = "  -3.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009061"
$ws.Range("E18").Value = "  -4.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.87"
$ws.Range("E20").Value = "  -4.80%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.926.19"
$ws.Range("E21").Value = "  -3.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.124"
$ws.Range("E22").Value = "  -3.96%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.03"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.110.30"
$ws.Range("E24").Value = "  -2.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.001"
$ws.Range("E25").Value = "  -4.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.17"
$ws.Range("E26").Value = "  -2.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.43"
$ws.Range("E27").Value = "  -3.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.473"
$ws.Range("E28").Value = "  -4.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.20"
$ws.Range("E29").Value = "  -3.87%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.716"
$ws.Range("E30").Value = "  -7.68%  "

# Row 31
$ws.Range("E31").Value = "  -3.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8223"
$ws.Range("E32").Value = "  -5.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.816"
$ws.Range("E33").Value = "  -5.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.179"
$ws.Range("E34").Value = "  -5.93%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.982"
$ws.Range("E35").Value = "  -2.20%  "

# Row 36
$ws.Range("E36").Value = "  +0.10%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.121"
$ws.Range("E37").Value = "  -2.86%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05501"
$ws.Range("E38").Value = "  -3.65%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01980"
$ws.Range("E39").Value = "  -3.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.968"
$ws.Range("E40").Value = "  -2.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5264"
$ws.Range("E41").Value = "  -4.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.041"
$ws.Range("E42").Value = "  -6.53%  "

# Row 43
$ws.Range("E43").Value = "  -2.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.780"
$ws.Range("E44").Value = "  -6.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06760"
$ws.Range("E45").Value = "  -2.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4900"
$ws.Range("E46").Value = "  -5.30%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.60"
$ws.Range("E47").Value = "  -6.08%  "

# Row 48
$ws.Range("E48").Value = "  -2.90%  "

# Row 50
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.882"
$ws.Range("E51").Value = "  -14.50%  "
